$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the extra "histórico" picture shapes (Image 3..6) from the drawing ---
$ws.Shapes.Item("Image 3").Delete()
$ws.Shapes.Item("Image 4").Delete()
$ws.Shapes.Item("Image 5").Delete()
$ws.Shapes.Item("Image 6").Delete()

# --- Clear the sample/demo data entered in the "histórico" table (row 18) ---
$ws.Range("A18").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""

# --- Remove the merged cells that belonged to the now-removed sample rows 19-22 ---
$ws.Range("B19:D19").UnMerge()
$ws.Range("E19:F19").UnMerge()
$ws.Range("J19:K19").UnMerge()
$ws.Range("B20:D20").UnMerge()
$ws.Range("E20:F20").UnMerge()
$ws.Range("J20:K20").UnMerge()
$ws.Range("B21:D21").UnMerge()
$ws.Range("E21:F21").UnMerge()
$ws.Range("J21:K21").UnMerge()
$ws.Range("B22:D22").UnMerge()
$ws.Range("E22:F22").UnMerge()
$ws.Range("J22:K22").UnMerge()

# --- Reset rows 19-22 back to their blank / unfilled (CRUD-ready) template state ---
$ws.Range("A19:K22").Clear()

$ws.Rows.Item(19).RowHeight = 12.75
$ws.Rows.Item(20).RowHeight = 12.75
$ws.Rows.Item(21).RowHeight = 12.75
$ws.Rows.Item(22).RowHeight = 15.75

# Re-apply the plain "unfilled row" style (same one used on row 16) to rows 19-21
$ws.Cells.Item(16,2).Copy()
foreach ($addr in @("A19","B19","C19","E19","F19","G19","H19","I19","J19","K19",
                     "A20","B20","C20","E20","F20","G20","H20","I20","J20","K20",
                     "A21","B21","C21","D21","E21","F21","G21","H21","I21","J21","K21")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
